$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 4.4620803418763
$ws.Cells.Item(2, 5).Value = 0.0487316900224034
$ws.Cells.Item(2, 6).Value = 9.17324141384519
$ws.Cells.Item(2, 7).Value = 1

$ws.Cells.Item(3, 4).Value = -1.94436936968067
$ws.Cells.Item(3, 5).Value = -13.528432402614
$ws.Cells.Item(3, 6).Value = 12.3223813947637
$ws.Cells.Item(3, 7).Value = 0

$ws.Cells.Item(4, 4).Value = -6.86483921876909
$ws.Cells.Item(4, 5).Value = -15.8793496859715
$ws.Cells.Item(4, 6).Value = 4.49968362782727
$ws.Cells.Item(4, 7).Value = 0

$ws.Cells.Item(5, 4).Value = -1.0734231285165
$ws.Cells.Item(5, 5).Value = -18.292886636882
$ws.Cells.Item(5, 6).Value = 18.9355016824991
$ws.Cells.Item(5, 7).Value = 0

$ws.Cells.Item(6, 4).Value = -45.8161278523735
$ws.Cells.Item(6, 5).Value = -52.4625948548407
$ws.Cells.Item(6, 6).Value = -38.0371930503741
$ws.Cells.Item(6, 7).Value = 1

$ws.Cells.Item(7, 4).Value = 20.310687173548
$ws.Cells.Item(7, 5).Value = 12.7696283939118
$ws.Cells.Item(7, 6).Value = 28.2935380333947
$ws.Cells.Item(7, 7).Value = 1

$ws.Cells.Item(8, 4).Value = 20.9120697187212
$ws.Cells.Item(8, 5).Value = 11.3063218609538
$ws.Cells.Item(8, 6).Value = 31.7074286495627
$ws.Cells.Item(8, 7).Value = 1

$ws.Cells.Item(9, 4).Value = 57.7048358712487
$ws.Cells.Item(9, 5).Value = 30.1121102994413
$ws.Cells.Item(9, 6).Value = 92.6300489366462
$ws.Cells.Item(9, 7).Value = 1

$ws.Cells.Item(10, 4).Value = 56.2263396253209
$ws.Cells.Item(10, 5).Value = 28.9858016511982
$ws.Cells.Item(10, 6).Value = 95.0920997190956
$ws.Cells.Item(10, 7).Value = 1

$ws.Cells.Item(11, 4).Value = -16.7221882581324
$ws.Cells.Item(11, 5).Value = -25.4740557857607
$ws.Cells.Item(11, 6).Value = -6.87793350549224
$ws.Cells.Item(11, 7).Value = 1

$ws.Cells.Item(12, 4).Value = 3.42926773014022
$ws.Cells.Item(12, 5).Value = -0.649730002167931
$ws.Cells.Item(12, 6).Value = 7.9168602605143
$ws.Cells.Item(12, 7).Value = 0

$ws.Cells.Item(13, 4).Value = 19.717269776304
$ws.Cells.Item(13, 5).Value = 7.64556013582691
$ws.Cells.Item(13, 6).Value = 34.1580900871777
$ws.Cells.Item(13, 7).Value = 1

$ws.Cells.Item(14, 4).Value = 119.941233146547
$ws.Cells.Item(14, 5).Value = 64.0966515107077
$ws.Cells.Item(14, 6).Value = 212.337419169173
$ws.Cells.Item(14, 7).Value = 1

$ws.Cells.Item(15, 4).Value = 205.336789972671
$ws.Cells.Item(15, 5).Value = 99.8184481009487
$ws.Cells.Item(15, 6).Value = 403.467035329428
$ws.Cells.Item(15, 7).Value = 1

$ws.Cells.Item(16, 4).Value = 146.866805185933
$ws.Cells.Item(16, 5).Value = 85.6710881283994
$ws.Cells.Item(16, 6).Value = 244.687752334718
$ws.Cells.Item(16, 7).Value = 1

$ws.Cells.Item(17, 4).Value = 4.04281806643453
$ws.Cells.Item(17, 5).Value = 0.854746669993673
$ws.Cells.Item(17, 6).Value = 7.53884418729428
$ws.Cells.Item(17, 7).Value = 1

$ws.Cells.Item(18, 4).Value = 21.7792733865914
$ws.Cells.Item(18, 5).Value = 14.0776701889269
$ws.Cells.Item(18, 6).Value = 31.3177937331251
$ws.Cells.Item(18, 7).Value = 1

$ws.Cells.Item(19, 4).Value = 75.4664393431562
$ws.Cells.Item(19, 5).Value = 52.6867297592668
$ws.Cells.Item(19, 6).Value = 106.762433359726
$ws.Cells.Item(19, 7).Value = 1

$ws.Cells.Item(20, 4).Value = 67.9517920871551
$ws.Cells.Item(20, 5).Value = 32.0786288686597
$ws.Cells.Item(20, 6).Value = 118.261046221184
$ws.Cells.Item(20, 7).Value = 1

$ws.Cells.Item(21, 4).Value = -43.5048466077246
$ws.Cells.Item(21, 5).Value = -48.730212751521
$ws.Cells.Item(21, 6).Value = -36.8311276706547
$ws.Cells.Item(21, 7).Value = 1

